$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Your account page" test-case block (rows 29-35) — fill in the Test Case ID
# column (B) for each existing row, clear the stray "Delete link" text that
# was left over in C34, and add a brand new row 35 test case for the
# "Your Account Page" redirect check (replacing the old stray "Buy Now" cell).
# ---------------------------------------------------------------------------
$ws.Range("B29").Value = "TC1"
$ws.Range("B30").Value = "TC2"
$ws.Range("B31").Value = "TC3"
$ws.Range("B32").Value = "TC4"
$ws.Range("B33").Value = "TC5"

$ws.Range("C34").Value = ""

$ws.Range("A35").Value = "Your Account Page "
$ws.Range("B35").Value = "TC1"
$ws.Range("C35").Value = "Verify User gets redirects to correct page"

# ---------------------------------------------------------------------------
# View state: scroll the sheet so row 19 is at the top and the active
# selection is C28 (matches where the author was working).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C28").Select()
